# Updates cryptos list (prices / 1h volume change) to match the latest
# scrape, and swaps the Aptos / PEPE rows (29 and 30) including their
# links and prices.
#
# NOTE: several "Price" values are plain decimal strings (e.g. "1.00",
# "0.540") that Excel's COM layer would otherwise silently coerce into
# numbers (stripping significant trailing zeros / formatting). To keep
# them as text - matching the original inline-string cells - those
# values are written with a leading apostrophe (forces text entry) and
# then the cell Style is reset back to "Normal" so no stray number
# formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Cell,
        [string]$Text,
        [bool]$Ambiguous = $false
    )
    if ($Ambiguous) {
        $ws.Range($Cell).Value = "'" + $Text
        $ws.Range($Cell).Style = "Normal"
    } else {
        $ws.Range($Cell).Value = $Text
    }
}

# Row 2 - Bitcoin
Set-TextValue "D2" "53.878.70"
Set-TextValue "E2" "  -0.81%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.250.44"
Set-TextValue "E3" "  -0.63%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.21%  "

# Row 5 - BNB
Set-TextValue "D5" "493.22" $true
Set-TextValue "E5" "  -0.57%  "

# Row 6 - Solana
Set-TextValue "D6" "127.97" $true
Set-TextValue "E6" "  +0.29%  "

# Row 7 - USDC
Set-TextValue "D7" "1.00" $true
Set-TextValue "E7" "  +0.02%  "

# Row 8 - XRP
Set-TextValue "D8" "0.522" $true
Set-TextValue "E8" "  -1.42%  "

# Row 9 - Dogecoin
Set-TextValue "E9" "  -0.43%  "

# Row 10 - TRON
Set-TextValue "E10" "  +0.75%  "

# Row 11 - Cardano
Set-TextValue "D11" "0.335" $true
Set-TextValue "E11" "  +3.11%  "

# Row 12 - Toncoin
Set-TextValue "D12" "4.72" $true
Set-TextValue "E12" "  +1.83%  "

# Row 13 - Wrapped liquid staked Ether 2.0
Set-TextValue "D13" "2.649.75"
Set-TextValue "E13" "  -1.57%  "

# Row 14 - Avalanche
Set-TextValue "D14" "22.58" $true
Set-TextValue "E14" "  +3.67%  "

# Row 15 - Wrapped BTC
Set-TextValue "D15" "53.863.44"
Set-TextValue "E15" "  -1.07%  "

# Row 16 - Shiba Inu
Set-TextValue "E16" "  -0.38%  "

# Row 17 - Wrapped Ether
Set-TextValue "D17" "2.259.34"
Set-TextValue "E17" "  -1.30%  "

# Row 18 - Chainlink
Set-TextValue "D18" "10.19" $true
Set-TextValue "E18" "  +1.24%  "

# Row 19 - Polkadot
Set-TextValue "E19" "  -0.16%  "

# Row 20 - Bitcoin Cash
Set-TextValue "D20" "299.59" $true
Set-TextValue "E20" "  -0.84%  "

# Row 21 - Uniswap
Set-TextValue "D21" "6.28" $true
Set-TextValue "E21" "  -2.88%  "

# Row 22 - Dai
Set-TextValue "E22" "  +0.11%  "

# Row 23 - Litecoin
Set-TextValue "D23" "60.57" $true
Set-TextValue "E23" "  -3.27%  "

# Row 24 - Binance-Peg BSC-USD
Set-TextValue "D24" "1.00" $true
Set-TextValue "E24" "  +0.33%  "

# Row 25 - Kaspa
Set-TextValue "E25" "  -1.95%  "

# Row 26 - Internet Computer (DFINITY)
Set-TextValue "E26" "  +2.49%  "

# Row 27 - Monero
Set-TextValue "D27" "170.18" $true
Set-TextValue "E27" "  +0.45%  "

# Row 28 - PancakeSwap
Set-TextValue "E28" "  -0.63%  "

# Row 29 - was Aptos, now PEPE
Set-TextValue "B29" "PEPE"
Set-TextValue "C29" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D29" "0.0₃0683"
Set-TextValue "E29" "  -0.65%  "

# Row 30 - was PEPE, now Aptos
Set-TextValue "B30" "Aptos"
Set-TextValue "C30" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D30" "5.90" $true
Set-TextValue "E30" "  +0.15%  "

# Row 31 - Fetch.AI
Set-TextValue "E31" "  +0.45%  "

# Row 32 - USDe
Set-TextValue "E32" "  -0.02%  "

# Row 33 - Ethereum Classic
Set-TextValue "D33" "17.68" $true
Set-TextValue "E33" "  +0.25%  "

# Row 34 - First Digital USD
Set-TextValue "D34" "0.998" $true
Set-TextValue "E34" "  +0.41%  "

# Row 35 - Sui Network
Set-TextValue "D35" "0.937" $true
Set-TextValue "E35" "  +8.20%  "

# Row 36 - Immutable X
Set-TextValue "D36" "1.18" $true
Set-TextValue "E36" "  -0.58%  "

# Row 37 - NEAR Protocol
Set-TextValue "D37" "3.68" $true
Set-TextValue "E37" "  -1.55%  "

# Row 38 - Polygon Ecosystem Token
Set-TextValue "D38" "0.369" $true

# Row 39 - Stacks
Set-TextValue "D39" "1.38" $true
Set-TextValue "E39" "  -1.54%  "

# Row 40 - Filecoin
Set-TextValue "D40" "3.34" $true
Set-TextValue "E40" "  -0.27%  "

# Row 41 - Aave
Set-TextValue "D41" "125.12" $true
Set-TextValue "E41" "  -2.20%  "

# Row 42 - Render Token
Set-TextValue "D42" "4.75" $true
Set-TextValue "E42" "  -1.23%  "

# Row 43 - Hedera
Set-TextValue "E43" "  +0.66%  "

# Row 44 - Stellar
Set-TextValue "D44" "0.0888" $true
Set-TextValue "E44" "  -0.42%  "

# Row 45 - Mantle
Set-TextValue "D45" "0.540" $true
Set-TextValue "E45" "  -0.75%  "

# Row 46 - Bittensor
Set-TextValue "D46" "237.98" $true
Set-TextValue "E46" "  -0.80%  "

# Row 47 - Polygon
Set-TextValue "E47" "  -1.46%  "

# Row 48 - VeChain
Set-TextValue "E48" "  -0.12%  "

# Row 49 - WhiteBIT Coin
Set-TextValue "E49" "  +0.26%  "

# Row 50 - Injective Protocol
Set-TextValue "D50" "16.03" $true
Set-TextValue "E50" "  -1.96%  "

# Row 51 - ZEEBU
Set-TextValue "E51" "  -1.10%  "
